$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.488.69'
$ws.Range('E2').Value = '  +0.38%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.105.82'
$ws.Range('E3').Value = '  +4.61%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '329.13'
$ws.Range('E5').Value = '  +1.22%  '
$ws.Range('E6').Value = '  +0.01%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5258'
$ws.Range('E7').Value = '  +2.63%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4385'
$ws.Range('E8').Value = '  +2.72%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.08869'
$ws.Range('E9').Value = '  +1.71%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '47.50'
$ws.Range('E10').Value = '  +9.65%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.163'
$ws.Range('E11').Value = '  +2.39%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '24.61'
$ws.Range('E12').Value = '  -0.13%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.104.14'
$ws.Range('E13').Value = '  +4.50%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.729'
$ws.Range('E14').Value = '  +1.97%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.763'
$ws.Range('E15').Value = '  +4.03%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '96.26'
$ws.Range('E16').Value = '  +2.11%  '
$ws.Range('E17').Value = '  +0.14%  '
$ws.Range('E18').Value = '  +1.30%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06636'
$ws.Range('E19').Value = '  +1.74%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '19.01'
$ws.Range('E20').Value = '  +0.58%  '
$ws.Range('E21').Value = '  +0.09%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.282'
$ws.Range('E22').Value = '  +1.19%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '30.572.91'
$ws.Range('E23').Value = '  +0.51%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.29'
$ws.Range('E24').Value = '  +3.81%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.351'
$ws.Range('E25').Value = '  +3.89%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.355.12'
$ws.Range('E26').Value = '  +4.71%  '
$ws.Range('E27').Value = '  -0.04%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.607'
$ws.Range('E28').Value = '  +7.22%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '161.83'
$ws.Range('E29').Value = '  -0.35%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '132.55'
$ws.Range('E30').Value = '  +1.10%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.208'
$ws.Range('E31').Value = '  +5.71%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.1075'
$ws.Range('E32').Value = '  +2.07%  '
$ws.Range('E33').Value = '  +22.38%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.208'
$ws.Range('E34').Value = '  +1.91%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.930'
$ws.Range('E35').Value = '  +2.65%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '10.21'
$ws.Range('E36').Value = '  +11.70%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02581'
$ws.Range('E37').Value = '  +2.22%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.490'
$ws.Range('E38').Value = '  +0.57%  '
$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06696'
$ws.Range('E39').Value = '  +0.34%  '
$ws.Range('B40').Value = 'Aptos'
$ws.Range('C40').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '12.71'
$ws.Range('E40').Value = '  +2.79%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.2288'
$ws.Range('E41').Value = '  +4.30%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.6851'
$ws.Range('E42').Value = '  +2.89%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.266'
$ws.Range('E43').Value = '  +2.23%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.001'
$ws.Range('E44').Value = '  +0.10%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '14.06'
$ws.Range('E45').Value = '  +3.04%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.6373'
$ws.Range('E46').Value = '  +3.34%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.211'
$ws.Range('E47').Value = '  +0.78%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.627'
$ws.Range('E48').Value = '  -1.15%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.252'
$ws.Range('E49').Value = '  -0.73%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.199'
$ws.Range('E50').Value = '  +8.25%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '82.28'
$ws.Range('E51').Value = '  +2.00%  '
